$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.032878
$ws.Range("H2").Value = 39.098634
$ws.Range("I2").Value = 0.02949184097968156
$ws.Range("J2").Value = 0.02949184097968156
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 1525.427510142314
$ws.Range("R2").Value = 13728.84759128083
$ws.Range("S2").Value = 0.009571180753587873
$ws.Range("T2").Value = 0.009571180753587873

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.032878
$ws.Range("H3").Value = 39.098634
$ws.Range("I3").Value = 0.02949184097968156
$ws.Range("J3").Value = 0.02949184097968156
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 1323.880233800779
$ws.Range("R3").Value = 11914.92210420701
$ws.Range("S3").Value = 0.008306587451426836
$ws.Range("T3").Value = 0.008306587451426834

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.032878
$ws.Range("H4").Value = 39.098634
$ws.Range("I4").Value = 0.02949184097968156
$ws.Range("J4").Value = 0.02949184097968156
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 1851.017818113035
$ws.Range("R4").Value = 16659.16036301732
$ws.Range("S4").Value = 0.01161407277466685
$ws.Range("T4").Value = 0.01161407277466685

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 395.9197996666667
$ws.Range("H5").Value = 1187.759399
$ws.Range("I5").Value = 0.8959190573622122
$ws.Range("J5").Value = 0.8959190573622122
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 46340.25993503255
$ws.Range("R5").Value = 417062.339415293
$ws.Range("S5").Value = 0.2907584929847395
$ws.Range("T5").Value = 0.2907584929847395

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 395.9197996666667
$ws.Range("H6").Value = 1187.759399
$ws.Range("I6").Value = 0.8959190573622122
$ws.Range("J6").Value = 0.8959190573622122
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 40217.54803114586
$ws.Range("R6").Value = 361957.9322803127
$ws.Range("S6").Value = 0.2523419953507245
$ws.Range("T6").Value = 0.2523419953507244

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 395.9197996666667
$ws.Range("H7").Value = 1187.759399
$ws.Range("I7").Value = 0.8959190573622122
$ws.Range("J7").Value = 0.8959190573622122
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 56231.21798015321
$ws.Range("R7").Value = 506080.9618213788
$ws.Range("S7").Value = 0.3528185690267482
$ws.Range("T7").Value = 0.3528185690267482

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 32.96202033333334
$ws.Range("H8").Value = 98.886061
$ws.Range("I8").Value = 0.07458910165810628
$ws.Range("J8").Value = 0.07458910165810628
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 3858.025265512114
$ws.Range("R8").Value = 34722.22738960903
$ws.Range("S8").Value = 0.0242068908044541
$ws.Range("T8").Value = 0.0242068908044541

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 32.96202033333334
$ws.Range("H9").Value = 98.886061
$ws.Range("I9").Value = 0.07458910165810628
$ws.Range("J9").Value = 0.07458910165810628
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 3348.283256042093
$ws.Range("R9").Value = 30134.54930437883
$ws.Range("S9").Value = 0.02100855271372469
$ws.Range("T9").Value = 0.02100855271372469

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.96202033333334
$ws.Range("H10").Value = 98.886061
$ws.Range("I10").Value = 0.07458910165810628
$ws.Range("J10").Value = 0.07458910165810628
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 4681.489917883384
$ws.Range("R10").Value = 42133.40926095045
$ws.Range("S10").Value = 0.02937365813992749
$ws.Range("T10").Value = 0.02937365813992749

